# Generate Report for Handoff
#
# The localization job has finished translation and the xliff files were
# (re)generated for handoff, so the status report is refreshed:
#   - Status moves from "In Translation" to "Ready for handoff" on every
#     sheet (the Overview roll-up as well as each per-locale sheet).
#   - The "Latest HO Xliff Generate Date"/"Latest Handoff Datetime"
#     timestamps are bumped to the new generation time.
#   - The Status column is widened so the longer "Ready for handoff"
#     text is no longer truncated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Width (in characters) to apply to the Status columns so that "Ready for
# handoff" fits without truncation, on the Overview sheet (columns E/F)
# and on each locale sheet (column C, "Status").
$statusColumnWidth = 16.33

# --- Overview sheet: zh-cn/de-de status + generate-date roll-up ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-28 16:59:02"

$wsOverview.Columns.Item(5).ColumnWidth = $statusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColumnWidth

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-28 16:58:56"

$wsZhCn.Columns.Item(3).ColumnWidth = $statusColumnWidth

# --- de-de sheet: Status + Latest Handoff Datetime ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-28 16:59:02"

$wsDeDe.Columns.Item(3).ColumnWidth = $statusColumnWidth
